$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look numeric (e.g. "0.999", "1.00") must be
# forced to Text format first, otherwise Excel auto-converts them to
# numbers on assignment and loses the original text formatting (e.g. trailing zeros).
$textFormatCells = @(
    "D4",
    "D5",
    "D6",
    "D7",
    "D10",
    "D12",
    "D15",
    "D17",
    "D19",
    "D20",
    "D21",
    "D23",
    "D24",
    "D25",
    "D26",
    "D27",
    "D28",
    "D29",
    "D30",
    "D31",
    "D32",
    "D33",
    "D34",
    "D35",
    "D36",
    "D37",
    "D38",
    "D39",
    "D40",
    "D41",
    "D42",
    "D43",
    "D44",
    "D45",
    "D47",
    "D48",
    "D49",
    "D51"
)
foreach ($ref in $textFormatCells) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range("D2").Value = "51.359.02"
$ws.Range("E2").Value = "  -1.83%  "
$ws.Range("D3").Value = "2.928.02"
$ws.Range("E3").Value = "  -2.43%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").Value = "374.10"
$ws.Range("E5").Value = "  +5.51%  "
$ws.Range("D6").Value = "104.70"
$ws.Range("E6").Value = "  -3.93%  "
$ws.Range("D7").Value = "0.547"
$ws.Range("E7").Value = "  -3.02%  "
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("E9").Value = "  -5.40%  "
$ws.Range("D10").Value = "37.27"
$ws.Range("E10").Value = "  -3.32%  "
$ws.Range("E11").Value = "  -0.36%  "
$ws.Range("D12").Value = "0.0841"
$ws.Range("E12").Value = "  -2.42%  "
$ws.Range("E13").Value = "  -4.88%  "
$ws.Range("D14").Value = "3.387.34"
$ws.Range("E14").Value = "  -2.45%  "
$ws.Range("D15").Value = "7.40"
$ws.Range("E15").Value = "  -4.99%  "
$ws.Range("D16").Value = "2.919.70"
$ws.Range("E16").Value = "  -2.71%  "
$ws.Range("D17").Value = "0.939"
$ws.Range("E17").Value = "  -8.61%  "
$ws.Range("D18").Value = "51.328.29"
$ws.Range("E18").Value = "  -1.99%  "
$ws.Range("D19").Value = "3.31"
$ws.Range("E19").Value = "  -6.11%  "
$ws.Range("D20").Value = "7.33"
$ws.Range("E20").Value = "  -3.21%  "
$ws.Range("D21").Value = "13.06"
$ws.Range("E21").Value = "  -4.67%  "
$ws.Range("E22").Value = "  -2.86%  "
$ws.Range("D23").Value = "68.71"
$ws.Range("E23").Value = "  -1.39%  "
$ws.Range("D24").Value = "261.40"
$ws.Range("E24").Value = "  -1.33%  "
$ws.Range("D25").Value = "2.70"
$ws.Range("E25").Value = "  -1.63%  "
$ws.Range("D26").Value = "0.171"
$ws.Range("E26").Value = "  -4.90%  "
$ws.Range("B27").Value = "Dai"
$ws.Range("C27").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  -0.04%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "25.91"
$ws.Range("E28").Value = "  -4.11%  "
$ws.Range("B29").Value = "Filecoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D29").Value = "7.17"
$ws.Range("E29").Value = "  -7.11%  "
$ws.Range("D30").Value = "6.86"
$ws.Range("E30").Value = "  +6.85%  "
$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D31").Value = "0.104"
$ws.Range("E31").Value = "  -4.40%  "
$ws.Range("B32").Value = "Cosmos"
$ws.Range("C32").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D32").Value = "9.91"
$ws.Range("E32").Value = "  -4.47%  "
$ws.Range("B33").Value = "Toncoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D33").Value = "2.12"
$ws.Range("E33").Value = "  -3.50%  "
$ws.Range("B34").Value = "InjectiveProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D34").Value = "34.72"
$ws.Range("E34").Value = "  -5.63%  "
$ws.Range("B35").Value = "OKB"
$ws.Range("C35").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D35").Value = "51.12"
$ws.Range("E35").Value = "  +0.28%  "
$ws.Range("B36").Value = "FirstDigitalUSD"
$ws.Range("C36").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").Value = "  +0.50%  "
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").Value = "0.0427"
$ws.Range("E37").Value = "  -4.39%  "
$ws.Range("B38").Value = "LidoDAOToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D38").Value = "3.02"
$ws.Range("E38").Value = "  -6.10%  "
$ws.Range("B39").Value = "Celestia"
$ws.Range("C39").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D39").Value = "17.21"
$ws.Range("E39").Value = "  -4.36%  "
$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D40").Value = "1.85"
$ws.Range("E40").Value = "  -6.45%  "
$ws.Range("D41").Value = "2.57"
$ws.Range("E41").Value = "  -5.10%  "
$ws.Range("B42").Value = "Stellar"
$ws.Range("C42").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D42").Value = "0.114"
$ws.Range("E42").Value = "  -3.92%  "
$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").Value = "22.00"
$ws.Range("E43").Value = "  -4.26%  "
$ws.Range("B44").Value = "Monero"
$ws.Range("C44").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D44").Value = "119.37"
$ws.Range("E44").Value = "  -3.64%  "
$ws.Range("B45").Value = "WEMIXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D45").Value = "2.09"
$ws.Range("E45").Value = "  -2.50%  "
$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").Value = "2.029.62"
$ws.Range("E46").Value = "  -4.61%  "
$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").Value = "2.31"
$ws.Range("E47").Value = "  -3.68%  "
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").Value = "3.21"
$ws.Range("E48").Value = "  -5.87%  "
$ws.Range("B49").Value = "TheGraph"
$ws.Range("C49").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D49").Value = "0.261"
$ws.Range("E49").Value = "  +5.13%  "
$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").Value = "3.216.87"
$ws.Range("E50").Value = "  -2.20%  "
$ws.Range("B51").Value = "BEAM"
$ws.Range("C51").Value = "https://coinranking.com/coin/cYYMfXF4u+beam-beam"
$ws.Range("D51").Value = "0.0324"
$ws.Range("E51").Value = "  -3.29%  "
